$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.248
$ws.Range("D5").Value = 0.374
$ws.Range("E5").Value = 0.414
$ws.Range("F5").Value = 0.462
$ws.Range("G5").Value = 0.492
$ws.Range("H5").Value = 0.518

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.248
$ws.Range("D7").Value = 0.374
$ws.Range("E7").Value = 0.414

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.241
$ws.Range("D8").Value = 0.47
$ws.Range("E8").Value = 0.513
$ws.Range("F8").Value = 0.56
$ws.Range("G8").Value = 0.588
$ws.Range("H8").Value = 0.622

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.281
$ws.Range("C9").Value = 0.395
$ws.Range("D9").Value = 0.505
$ws.Range("E9").Value = 0.546
$ws.Range("G9").Value = 0.56
$ws.Range("H9").Value = 0.596
